# Insert two new data rows (187 and 188) into the "Haba" sheet, pushing
# the existing rows 187-283 down to 189-285, and populate the two new
# rows with their own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 187, shifting
# everything from (old) row 187 downward by two rows.
$ws.Rows.Item(187).Insert()
$ws.Rows.Item(188).Insert()

# --- New row 187 ---
$ws.Cells.Item(187, 1).Value = 6
$ws.Cells.Item(187, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(187, 3).Value = "Metropolitana"
$ws.Cells.Item(187, 4).Value = 44784
$ws.Cells.Item(187, 5).Value = 13
$ws.Cells.Item(187, 6).Value = 100112026
$ws.Cells.Item(187, 7).Value = "Haba"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 500
$ws.Cells.Item(187, 11).Value = 14000
$ws.Cells.Item(187, 12).Value = 15000
$ws.Cells.Item(187, 13).Value = 14540
$ws.Cells.Item(187, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(187, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(187, 16).Value = 582
$ws.Cells.Item(187, 17).Value = 25
$ws.Cells.Item(187, 18).Value = "Hortaliza"

# --- New row 188 ---
$ws.Cells.Item(188, 1).Value = 6
$ws.Cells.Item(188, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(188, 3).Value = "Metropolitana"
$ws.Cells.Item(188, 4).Value = 44784
$ws.Cells.Item(188, 5).Value = 13
$ws.Cells.Item(188, 6).Value = 100112026
$ws.Cells.Item(188, 7).Value = "Haba"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Segunda"
$ws.Cells.Item(188, 10).Value = 180
$ws.Cells.Item(188, 11).Value = 12000
$ws.Cells.Item(188, 12).Value = 12000
$ws.Cells.Item(188, 13).Value = 12000
$ws.Cells.Item(188, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(188, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(188, 16).Value = 480
$ws.Cells.Item(188, 17).Value = 25
$ws.Cells.Item(188, 18).Value = "Hortaliza"
